# Auto-applies numeric cell updates produced by the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 425.66666
$ws.Range("J4").Value = 440
$ws.Range("L4").Value = 440
$ws.Range("N4").Value = -668
$ws.Range("H10").Value = 1000
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H33").Value = 218.33333
$ws.Range("I33").Value = 249.61539
$ws.Range("J33").Value = 137
$ws.Range("K33").Value = 249.61539
$ws.Range("L33").Value = 137
$ws.Range("M33").Value = -20.61538999999999
$ws.Range("N33").Value = -595
$ws.Range("H58").Value = 2337.2222
$ws.Range("I58").Value = 321.375
$ws.Range("K58").Value = 964.125
$ws.Range("M58").Value = -814.125
$ws.Range("H62").Value = 5230
$ws.Range("I62").Value = 5212.4287
$ws.Range("J62").Value = 5353
$ws.Range("K62").Value = 5212.4287
$ws.Range("L62").Value = 5353
$ws.Range("M62").Value = -4588.4287
$ws.Range("N62").Value = -6601
$ws.Range("H65").Value = 5230
$ws.Range("I65").Value = 5212.4287
$ws.Range("J65").Value = 5353
$ws.Range("K65").Value = 26062.1435
$ws.Range("L65").Value = 26765
$ws.Range("M65").Value = -22942.1435
$ws.Range("N65").Value = -33005
$ws.Range("H98").Value = 1663.5
$ws.Range("I98").Value = 1610.6666
$ws.Range("K98").Value = 1610.6666
$ws.Range("M98").Value = -112.6666
$ws.Range("H100").Value = 1724.2307
$ws.Range("I100").Value = 1256.5555
$ws.Range("K100").Value = 1256.5555
$ws.Range("M100").Value = -715.5554999999999
$ws.Range("H122").Value = 1663.5
$ws.Range("I122").Value = 1610.6666
$ws.Range("K122").Value = 4831.9998
$ws.Range("M122").Value = -2381.9998
$ws.Range("H134").Value = 69999.5
$ws.Range("J134").Value = 69999.5
$ws.Range("L134").Value = 69999.5
$ws.Range("N134").Value = -80139.5
$ws.Range("H137").Value = 13374.333
$ws.Range("I137").Value = 17011.1
$ws.Range("K137").Value = 51033.3
$ws.Range("M137").Value = -48483.3
$ws.Range("H138").Value = 33201.97
$ws.Range("I138").Value = 2547.4119
$ws.Range("J138").Value = 65772.44
$ws.Range("K138").Value = 7642.2357
$ws.Range("L138").Value = 197317.32
$ws.Range("M138").Value = -2502.2357
$ws.Range("N138").Value = -207597.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23874.355
$ws.Range("I32").Value = 23874.355
$ws.Range("K32").Value = 23874.355
$ws.Range("M32").Value = -23587.355
$ws.Range("H45").Value = 3607.4666
$ws.Range("J45").Value = 4760.625
$ws.Range("L45").Value = 4760.625
$ws.Range("N45").Value = -5514.625
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 1621.6613
$ws.Range("I132").Value = 1230.6945
$ws.Range("K132").Value = 3692.0835
$ws.Range("M132").Value = -1162.0835

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1388.7368
$ws.Range("I134").Value = 974.3674
$ws.Range("K134").Value = 2923.1022
$ws.Range("M134").Value = -388.1021999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2298.75
$ws.Range("I16").Value = 1546.1818
$ws.Range("K16").Value = 1546.1818
$ws.Range("M16").Value = -1259.1818
$ws.Range("H31").Value = 6252825.5
$ws.Range("I31").Value = 20001000
$ws.Range("J31").Value = 3655.818
$ws.Range("K31").Value = 20001000
$ws.Range("L31").Value = 3655.818
$ws.Range("M31").Value = -20000705
$ws.Range("N31").Value = -4245.818
$ws.Range("H34").Value = 6252825.5
$ws.Range("I34").Value = 20001000
$ws.Range("J34").Value = 3655.818
$ws.Range("K34").Value = 20001000
$ws.Range("L34").Value = 3655.818
$ws.Range("M34").Value = -20000798
$ws.Range("N34").Value = -4059.818
$ws.Range("H99").Value = 11712
$ws.Range("I99").Value = 9493.75
$ws.Range("K99").Value = 9493.75
$ws.Range("M99").Value = -7995.75
$ws.Range("H107").Value = 596.85297
$ws.Range("I107").Value = 504
$ws.Range("J107").Value = 819.7
$ws.Range("K107").Value = 504
$ws.Range("L107").Value = 819.7
$ws.Range("M107").Value = 1416
$ws.Range("N107").Value = -4659.7
$ws.Range("H113").Value = 2298.75
$ws.Range("I113").Value = 1546.1818
$ws.Range("K113").Value = 1546.1818
$ws.Range("M113").Value = 623.8181999999999
$ws.Range("H126").Value = 11712
$ws.Range("I126").Value = 9493.75
$ws.Range("K126").Value = 28481.25
$ws.Range("M126").Value = -26011.25
$ws.Range("H134").Value = 1327.64
$ws.Range("I134").Value = 1075.0244
$ws.Range("K134").Value = 3225.0732
$ws.Range("M134").Value = -690.0731999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 250000200
$ws.Range("J9").Value = 333333500
$ws.Range("L9").Value = 1000000500
$ws.Range("N9").Value = -1000000948
$ws.Range("H39").Value = 20000
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 60000
$ws.Range("N39").Value = -60588
$ws.Range("H117").Value = 2093.5
$ws.Range("I117").Value = 3499.8
$ws.Range("J117").Value = 1312.2222
$ws.Range("K117").Value = 10499.4
$ws.Range("L117").Value = 3936.6666
$ws.Range("M117").Value = -7057.400000000001
$ws.Range("N117").Value = -10820.6666
$ws.Range("H121").Value = 66501.62
$ws.Range("I121").Value = 152970
$ws.Range("J121").Value = 23267.428
$ws.Range("K121").Value = 458910
$ws.Range("L121").Value = 69802.284
$ws.Range("M121").Value = -457600
$ws.Range("N121").Value = -72422.284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15144.784
$ws.Range("I102").Value = 17605.871
$ws.Range("J102").Value = 2429.1667
$ws.Range("K102").Value = 17605.871
$ws.Range("L102").Value = 2429.1667
$ws.Range("M102").Value = -15983.871
$ws.Range("N102").Value = -5673.1667
$ws.Range("H122").Value = 3300.7273
$ws.Range("J122").Value = 4736.5
$ws.Range("L122").Value = 14209.5
$ws.Range("N122").Value = -19109.5
$ws.Range("H126").Value = 3279.1738
$ws.Range("I126").Value = 2592
$ws.Range("J126").Value = 4849.857
$ws.Range("K126").Value = 7776
$ws.Range("L126").Value = 14549.571
$ws.Range("M126").Value = -5306
$ws.Range("N126").Value = -19489.571
$ws.Range("H132").Value = 2126.0571
$ws.Range("J132").Value = 3998.889
$ws.Range("L132").Value = 11996.667
$ws.Range("N132").Value = -17056.667
$ws.Range("H133").Value = 93653.8
$ws.Range("J133").Value = 95390
$ws.Range("L133").Value = 95390
$ws.Range("N133").Value = -105510

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2372.0908
$ws.Range("I40").Value = 2109.3
$ws.Range("K40").Value = 2109.3
$ws.Range("M40").Value = -1973.3
$ws.Range("H45").Value = 90020.5
$ws.Range("I45").Value = 100041
$ws.Range("K45").Value = 100041
$ws.Range("M45").Value = -99634
$ws.Range("H132").Value = 2582.3635
$ws.Range("I132").Value = 1279.6
$ws.Range("K132").Value = 3838.8
$ws.Range("M132").Value = -1308.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 39118.49
$ws.Range("I122").Value = 46832.75
$ws.Range("J122").Value = 3853.2856
$ws.Range("K122").Value = 140498.25
$ws.Range("L122").Value = 11559.8568
$ws.Range("M122").Value = -138048.25
$ws.Range("N122").Value = -16459.8568
$ws.Range("H132").Value = 1632.1111
$ws.Range("I132").Value = 1182.5264
$ws.Range("K132").Value = 3547.5792
$ws.Range("M132").Value = -1017.5792
$ws.Range("H136").Value = 18120.629
$ws.Range("I136").Value = 19700.773
$ws.Range("K136").Value = 59102.319
$ws.Range("M136").Value = -56552.319
